# Rename the header row to lowercase, underscore-separated, R-friendly names.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "distance_in_miles"
$ws.Range("C1").Value = "gasoline_in_gallons"
$ws.Range("D1").Value = "comments"
